# Grupo 2D - Viajar 360: add the "viernes 05/02" meeting column (F)
# to the Scrum Master report sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column F: "viernes 05/02" meeting -------------------------------

$ws.Range("F5").Value  = "viernes 05/02"
$ws.Range("F6").Value  = "Ausente"
$ws.Range("F7").Value  = "Participó"
$ws.Range("F8").Value  = "Participó editando"
$ws.Range("F9").Value  = "Ausente"
$ws.Range("F10").Value = "Ausente"
$ws.Range("F11").Value = "Participó editando"
$ws.Range("F13").Value = "Realizamos diagrama de clases,`n Falta pulir métodos. La reunión fue grabada para los compañeros ausentes"

# --- Formatting: borders for the whole new column (F5:F13) ---------------

$newCol = $ws.Range("F5:F13")
$newCol.Borders.Color = 0
$newCol.Borders.LineStyle = 1
$newCol.Borders.Weight = 2

# Wrap text on the long note in F13, same as the rest of row 13
$ws.Range("F13").WrapText = $true

# Row 13 grows to fit the new, longer note text
$ws.Rows.Item(13).RowHeight = 57

# Column F width, matching the other data columns
$ws.Columns.Item(6).ColumnWidth = 27.59

# Match the author's last selection when they saved the file
$ws.Range("D16").Select()

Write-Host "Added viernes 05/02 column"
